# Updated cryptos list on Sat May 13 10:44:13 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# each coin row (rows 2-51) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '26.940.53'; E = '  +1.88%  ' },
    @{ Row = 3; D = '1.814.65'; E = '  +2.39%  ' },
    @{ Row = 4; D = '1.006'; E = '  +0.33%  ' },
    @{ Row = 5; D = '312.16'; E = '  +1.80%  ' },
    @{ Row = 6; D = '1.005'; E = '  +0.19%  ' },
    @{ Row = 7; D = $null; E = '  -0.06%  ' },
    @{ Row = 8; D = '0.3668'; E = '  -0.18%  ' },
    @{ Row = 9; D = '0.07256'; E = '  +0.17%  ' },
    @{ Row = 10; D = '2.145.31'; E = '  +20.84%  ' },
    @{ Row = 11; D = '0.8634'; E = '  +1.66%  ' },
    @{ Row = 12; D = '21.27'; E = '  +4.61%  ' },
    @{ Row = 13; D = '5.406'; E = '  +3.00%  ' },
    @{ Row = 14; D = '6.601'; E = '  +2.53%  ' },
    @{ Row = 15; D = '0.06972'; E = '  +2.23%  ' },
    @{ Row = 16; D = '80.99'; E = '  +1.68%  ' },
    @{ Row = 17; D = $null; E = '  +0.53%  ' },
    @{ Row = 18; D = '0.000008867'; E = '  +2.34%  ' },
    @{ Row = 19; D = '1.005'; E = '  +0.25%  ' },
    @{ Row = 20; D = '15.24'; E = '  +1.47%  ' },
    @{ Row = 21; D = '26.985.42'; E = '  +2.05%  ' },
    @{ Row = 22; D = '5.181'; E = '  +1.51%  ' },
    @{ Row = 23; D = $null; E = '  -2.78%  ' },
    @{ Row = 24; D = '2.354.51'; E = '  +18.15%  ' },
    @{ Row = 25; D = '153.75'; E = '  +0.85%  ' },
    @{ Row = 26; D = '1.883'; E = '  +1.71%  ' },
    @{ Row = 27; D = '18.29'; E = '  +0.72%  ' },
    @{ Row = 28; D = '5.222'; E = '  +2.49%  ' },
    @{ Row = 29; D = '1.901'; E = '  +10.76%  ' },
    @{ Row = 30; D = '114.62'; E = '  -0.20%  ' },
    @{ Row = 31; D = $null; E = '  +0.19%  ' },
    @{ Row = 32; D = $null; E = '  +6.52%  ' },
    @{ Row = 33; D = '0.7470'; E = '  +3.08%  ' },
    @{ Row = 34; D = '4.418'; E = '  +1.84%  ' },
    @{ Row = 35; D = '2.808'; E = '  +1.84%  ' },
    @{ Row = 36; D = '1.005'; E = '  +0.23%  ' },
    @{ Row = 37; D = '1.130'; E = '  +4.90%  ' },
    @{ Row = 38; D = '0.05211'; E = '  +1.04%  ' },
    @{ Row = 39; D = '0.01922'; E = '  +1.33%  ' },
    @{ Row = 40; D = '0.5098'; E = '  +3.28%  ' },
    @{ Row = 41; D = '0.1653'; E = '  +2.79%  ' },
    @{ Row = 42; D = '2.737'; E = '  +8.32%  ' },
    @{ Row = 43; D = '6.462'; E = '  +3.96%  ' },
    @{ Row = 44; D = '8.339'; E = '  +3.45%  ' },
    @{ Row = 45; D = '106.72'; E = '  +1.68%  ' },
    @{ Row = 46; D = '10.42'; E = '  +2.81%  ' },
    @{ Row = 47; D = '1.004'; E = '  +0.18%  ' },
    @{ Row = 48; D = '0.4563'; E = '  +1.34%  ' },
    @{ Row = 49; D = '1.643'; E = '  +3.81%  ' },
    @{ Row = 50; D = '0.06214'; E = '  +0.20%  ' },
    @{ Row = 51; D = '1.841'; E = '  +5.44%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)   # column D = Price

        if ($u.D -match '^-?\d+(\.\d+)?$') {
            # Looks like a plain number (e.g. "1.006"/"21.27") - force text
            # storage so Excel doesn't coerce it into a numeric value; the
            # source data is always stored as text in this sheet.
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }

    $ws.Cells.Item($row, 5).Value = $u.E   # column E = Volume(1h)
}
